# Weekly update: insert two new rows of data (row 84 and 85) above the
# existing block, shifting the rest of the table down by two rows.
# This mirrors the commit "Fruta / hortaliza, semanal" which adds a new
# week's worth of price entries to the top of the existing data block and
# pushes the older rows down (sheet dimension grows from R114 to R116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 84:85 - this shifts existing rows 84-114 down to
# rows 86-116 (carrying their values/styles with them), matching the diff.
$ws.Rows("84:85").Insert()

# New row 84 - Magnum
$ws.Range("A84").Value = 2
$ws.Range("B84").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44524
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = 100112031
$ws.Range("G84").Value = "Poroto verde"
$ws.Range("H84").Value = "Magnum"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 340
$ws.Range("K84").Value = 16000
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = 17000
$ws.Range("N84").Value = "$/malla 25 kilos"
$ws.Range("O84").Value = "Provincia de Limarí"
$ws.Range("P84").Value = 680
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"

# New row 85 - Sin especificar
$ws.Range("A85").Value = 2
$ws.Range("B85").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 44524
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = 100112031
$ws.Range("G85").Value = "Poroto verde"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 240
$ws.Range("K85").Value = 35000
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = 37500
$ws.Range("N85").Value = "$/malla 25 kilos"
$ws.Range("O85").Value = "Provincia de Limarí"
$ws.Range("P85").Value = 1500
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
